# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns
# with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.719.30"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "1.638.83"
$ws.Range("E3").Value = "  -0.49%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.07"
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.08"
$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Value = "1.866.50"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "1.653.76"
$ws.Range("E13").Value = "  +0.12%  "

$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.48"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("D17").Value = "26.693.99"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("E18").Value = "  -2.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "211.37"
$ws.Range("E19").Value = "  -3.18%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("E22").Value = "  -1.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("E23").Value = "  -4.84%  "

$ws.Range("E24").Value = "  -2.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.84"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("E30").Value = "  -3.21%  "

$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("E33").Value = "  -0.61%  "

$ws.Range("D34").Value = "1.268.94"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0174"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("E42").Value = "  -2.84%  "

$ws.Range("D43").Value = "1.777.85"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("E44").Value = "  -3.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.30"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.12"
$ws.Range("E46").Value = "  +0.90%  "

$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("E48").Value = "  +0.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.54"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("E51").Value = "  -0.26%  "
